$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(3.182878228561681, 1.65323645889881, 3.082599426703578, 0.4998867070740569, 0, 8.418600821238126),
    @(1.505614041169197, 1.65323645889881, 0.1529057820181812, 0.4998867070740569, 1, 3.811642989160245),
    @(0.3464964993005633, 1.65323645889881, 0.7127328510149897, 6.48142807727062, 0, 9.193893886484982),
    @(0.1554434735375247, 0.3375848360084654, 0.1529057820181812, 0.4998867070740569, 1, 1.145820798638228),
    @(1.505614041169197, 1.65323645889881, 0.1529057820181812, 0.4998867070740569, 0, 3.811642989160245),
    @(0.7287194209349384, 0.3375848360084654, 0.7127328510149897, 0.4998867070740569, 1, 2.27892381503245)
)

$cols = @("B", "C", "D", "E", "F", "G")

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 2
    for ($j = 0; $j -lt $cols.Length; $j++) {
        $ws.Range("$($cols[$j])$row").Value = $data[$i][$j]
    }
}
